# Updated cryptos list on Fri Dec 29 11:47:10 UTC 2023 with GitHub Actions
#
# Applies the latest price/volume refresh to the cryptos table on Sheet1.
# Column D ("Price") holds numeric-looking text (e.g. "319.09", "1.00") that
# must stay plain text (matching the source feed's inline-string cells), so
# values are entered with a leading apostrophe to force text entry, then the
# cell style is reset to "Normal" so no stray number-format style is left
# behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "42.957.07"
$ws.Range("E2").Value = "  -0.57%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.371.37"
$ws.Range("E3").Value = "  -1.35%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
Set-TextValue "D5" "319.09"
$ws.Range("E5").Value = "  -4.50%  "

# Row 6 - Solana
Set-TextValue "D6" "108.17"
$ws.Range("E6").Value = "  +2.14%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -2.41%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -4.15%  "

# Row 10 - Avalanche
Set-TextValue "D10" "41.98"
$ws.Range("E10").Value = "  -0.93%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -1.21%  "

# Row 12 - Polkadot
$ws.Range("E12").Value = "  -1.83%  "

# Row 13 - Polygon
$ws.Range("E13").Value = "  -4.64%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.02%  "

# Row 15 - Chainlink
$ws.Range("E15").Value = "  -6.17%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue "D16" "2.731.95"
$ws.Range("E16").Value = "  -1.08%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.386.53"
$ws.Range("E17").Value = "  -0.36%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "42.936.73"
$ws.Range("E18").Value = "  -0.64%  "

# Row 19 - Uniswap
Set-TextValue "D19" "7.68"
$ws.Range("E19").Value = "  -0.32%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -2.37%  "

# Row 21 - Litecoin
Set-TextValue "D21" "76.29"
$ws.Range("E21").Value = "  -1.23%  "

# Row 22 - PancakeSwap
Set-TextValue "D22" "3.71"
$ws.Range("E22").Value = "  -3.42%  "

# Row 23 - BitcoinCash
Set-TextValue "D23" "258.12"
$ws.Range("E23").Value = "  -6.85%  "

# Row 24 - ImmutableX
Set-TextValue "D24" "2.35"
$ws.Range("E24").Value = "  -3.28%  "

# Row 25 - InternetComputer(DFINITY)
Set-TextValue "D25" "9.45"
$ws.Range("E25").Value = "  -4.22%  "

# Row 26 - Dai
Set-TextValue "D26" "1.00"
$ws.Range("E26").Value = "  +0.14%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  -3.81%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "23.09"
$ws.Range("E28").Value = "  -0.85%  "

# Row 29 - Toncoin
Set-TextValue "D29" "2.26"
$ws.Range("E29").Value = "  +2.54%  "

# Row 30 / Row 31 - Monero and InjectiveProtocol swap rank order
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D30" "36.93"
$ws.Range("E30").Value = "  -0.96%  "

$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D31" "171.28"
$ws.Range("E31").Value = "  -2.60%  "

# Row 32 - Hedera
$ws.Range("E32").Value = "  -4.33%  "

# Row 33 - Filecoin
Set-TextValue "D33" "6.07"
$ws.Range("E33").Value = "  -0.73%  "

# Row 34 - WEMIXToken
$ws.Range("E34").Value = "  -6.47%  "

# Row 35 - Kaspa
Set-TextValue "D35" "0.122"
$ws.Range("E35").Value = "  +11.91%  "

# Row 36 - Stellar
$ws.Range("E36").Value = "  -3.25%  "

# Row 37 - RenderToken
Set-TextValue "D37" "4.74"
$ws.Range("E37").Value = "  -3.17%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -0.31%  "

# Row 39 - NEARProtocol
Set-TextValue "D39" "3.91"
$ws.Range("E39").Value = "  -4.97%  "

# Row 40 - LidoDAOToken
$ws.Range("E40").Value = "  -4.38%  "

# Row 41 - ARBITRUM
Set-TextValue "D41" "1.55"
$ws.Range("E41").Value = "  -1.30%  "

# Row 42 - Algorand
$ws.Range("E42").Value = "  +2.78%  "

# Row 43 - MultiversX
Set-TextValue "D43" "71.80"
$ws.Range("E43").Value = "  +1.90%  "

# Row 44 - FirstDigitalUSD
$ws.Range("E44").Value = "  -0.04%  "

# Row 45 - Celestia
$ws.Range("E45").Value = "  +0.59%  "

# Row 46 - BitcoinSV
Set-TextValue "D46" "90.55"
$ws.Range("E46").Value = "  -1.65%  "

# Row 47 - Aave
$ws.Range("E47").Value = "  -6.77%  "

# Row 48 - THORChain
$ws.Range("E48").Value = "  -0.46%  "

# Row 49 - FraxShare
$ws.Range("E49").Value = "  -0.36%  "

# Row 50 - ordi
Set-TextValue "D50" "77.49"
$ws.Range("E50").Value = "  +7.82%  "

# Row 51 - TrustWalletToken
$ws.Range("E51").Value = "  -1.92%  "
